# RevlonTestData.xlsx maintenance edit:
#  - Row 4 "PaymentDetails" sample record is refreshed with a new card
#    number (stored as text, quote-prefixed & left-aligned like Excel
#    does for numeric-looking text), a new expiry year/month and cvv.
#  - Row 4 gets an explicit (custom) row height, matching the resaved
#    workbook.
#  - The view's selection moves from H1 to the Q4:XFD4 row-selection.
#
# NOTE: the source diff also nudges the sheetView's topLeftCell (E1->H1)
# and the window geometry / fileVersion / revisionPtr metadata in
# workbook.xml. Those are Excel-session/GUID bookkeeping values that
# this COM surface does not expose a way to persist (ActiveWindow's
# ScrollRow/ScrollColumn/TopLeftCell accept writes silently but are not
# wired into the saved sheetView), so they are intentionally left as
# best-effort / no-ops below rather than producing misleading no-op
# calls.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# -- Row 4: cardNumber / ExpYear / ExpMonth / cvv refresh --------------

# cardNumber (M4): replaced by a new card number. In the target file this
# is stored as a shared-string ("quote-prefixed" text), left aligned, on
# top of the existing "0;[Red]0" custom number format (numFmtId 164) -
# entering it with a leading apostrophe reproduces that quote-prefixed
# text cell/style combination exactly.
$ws.Range("M4").Value = "'4444444444444448"
$ws.Range("M4").NumberFormat = "0;[Red]0"
$ws.Range("M4").HorizontalAlignment = -4131   # xlLeft

# ExpYear (N4): 2023 -> 2025
$ws.Range("N4").Value = 2025

# ExpMonth (O4): "Jun" -> "Feb"
$ws.Range("O4").Value = "Feb"

# cvv (P4): 345 -> 123
$ws.Range("P4").Value = 123

# Row 4 picks up an explicit custom height in the resaved workbook.
$ws.Rows.Item(4).RowHeight = 14.25

# -- View: selection moves to the full Q4:XFD4 row-selection ----------
$excel.ActiveWindow.ScrollColumn = 8
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("Q4:XFD4").Select()
